$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header labels for SVR parameters
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# Add the corresponding values
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 5

# Update selection to match saved view state
$ws.Range("J7").Select()
